$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated s_vals data (regen to filter save games)
$data = @{
    2 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245)
    3 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 6.741336633845642)
    4 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 1, 28.30127388105354)
    5 = @(0.3464964993005633, 9.226618575922256, 0.7127328510149897, 246.9852506941017, 0, 257.2710986203395)
    6 = @(0.1554434735375247, 1.65323645889881, 3.082599426703578, 6.48142807727062, 1, 11.37270743641053)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]  # B: TB
    $ws.Cells.Item($row, 3).Value = $values[1]  # C: d2S
    $ws.Cells.Item($row, 4).Value = $values[2]  # D: K
    $ws.Cells.Item($row, 5).Value = $values[3]  # E: IP
    $ws.Cells.Item($row, 6).Value = $values[4]  # F: Win
    $ws.Cells.Item($row, 7).Value = $values[5]  # G: sum
}
